$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all updated crypto values from the Aug 23 2024 GitHub Actions refresh.
# NumberFormat is set to Text ("@") before assignment so that numeric-looking
# strings (e.g. "1.00", "351.00", "20.50") are preserved exactly as text,
# matching the original inline-string cell contents instead of being coerced
# into trimmed numeric values by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.018.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.53%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.651.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.62%  '
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.04'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.46%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.25'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.26%  '
# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.17%  '
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.20%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.60'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.52%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.51%  '
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.24%  '
# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.08%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.142.63'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.14%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.88'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +11.25%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '61.034.76'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.60%  '
# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.31%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.669.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +8.22%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.70'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.11%  '
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.28%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '351.00'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.71%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.95'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.07%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.20%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.536'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.14%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.06'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.21%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.33%  '
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.01%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.16'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.00%  '
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.32%  '
# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.86%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.86'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +7.33%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.18%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '166.23'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.63%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.93'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.94%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.08'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +9.77%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.49'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.90%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.32'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.13%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.66'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.39%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '333.47'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +12.08%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.02'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.14%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.67'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.06%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.884'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.79%  '
# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.50'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.84%  '
# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.19'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.68%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '134.68'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.08%  '
# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0565'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.29%  '
# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.100'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.59%  '
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.96%  '
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.44%  '
# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.50'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.92%  '
# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.28%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.094.91'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.26%  '
